$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated data values (rows 2-53) ---
$row = @(-0.41785794355,-0.49401922234999995,-0.19453355088999999,0.18654143340000001,0.021261309300000003,-0.42325712641000002,-0.2717334937,-0.35862523774999999,-0.24163961393,0.11873503835,-0.11309856330000002,-0.31019287179499999,0.062167637070000017,0.208774914895)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(2, $i + 1).Value = $row[$i]
}
$row = @(0.13939963784999998,-0.090831402283333307,0.32564868923333329,0.26241120803333334,0.33007849423333335,0.095779555849999976,0.21429990949333333,-0.11760521124500001,-0.2455154556166666,-0.23538284410000004,0.37938201027000007,0.11883656286499999,-0.18339116655000004,1.3342409641849999)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(3, $i + 1).Value = $row[$i]
}
$row = @(-0.29426727150500004,-1.0087204467233333,-1.5376417652333334,-0.72671838121666665,-0.20290285209999998,-0.25132171556666666,-0.73760390350000005,-0.5930549971600001,-0.21227342002666666,-0.54607445902999996,-1.0290494167066666,0.022644864883333365,0.0076014471766666092,-0.22359010002333329)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(4, $i + 1).Value = $row[$i]
}
$row = @(-0.11519638369999999,-0.31523474541999996,-0.27026040723666672,0.55912296753333346,0.47098399390000006,0.030480387693333319,-0.33457160631666671,0.093453615359999978,0.98111972717999996,0.044964584400000021,0.0027975351566666279,0.30725543822000001,0.48600499934999997,1.3892934171833333)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(5, $i + 1).Value = $row[$i]
}
$row = @(-0.39760818776666668,-0.085617780935000026,0.11629278763333332,-0.86176866185000001,-1.0071449895666666,-0.58090549902999999,-0.71774609309333326,-0.73527436265000001,-0.22052423480333333,-1.0257276454333333,-0.01477055459999993,-0.47787223948333335,-0.98851786029666666,0.10766542534500001)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(6, $i + 1).Value = $row[$i]
}
$row = @(-1.1852303042400001,-0.82479976059999993,-0.80075912506999991,-0.35709857322999999,-0.48974067204833338,-0.33351493890000006,0.061043278596666639,0.16296812915333334,-0.88709962735983328,0.42375952239999998,-0.87879634776666671,-1.2977798774,-0.91777295893333344,-0.56612170009999996)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(7, $i + 1).Value = $row[$i]
}
$row = @(-0.1068290247,-0.48640855442999997,-0.30467925572333332,-0.55524245958333329,-0.62687483744999994,-0.42305760711666673,0.53381399283333331,0.094990716946666665,-0.43565054898333333,0.52679305650000008,-0.96124580640000001,-0.48445321233333338,-0.05497227799999993,0.037746142966666674)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(8, $i + 1).Value = $row[$i]
}
$row = @(0.49678745262999996,0.34229507080333332,0.43247000026666671,0.76719099171333338,1.318277578,0.30158052084600007,-0.57092418394666666,-0.78174125558333329,0.5440688879433333,-0.084649662199999698,-0.4297574036516667,0.13947828806333329,-0.25769499130000001,1.1249797051133332)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(9, $i + 1).Value = $row[$i]
}
$row = @(0.20083170721666666,-0.0078327315666666952,0.058584598593333315,0.23564532483,0.00096246346666661875,-0.061839679666666647,0.15887030295000001,0.14051754377000003,-0.15285195817133335,0.071607939549999985,0.20432310075999999,0.55551785026999989,-0.31139521811666665,-0.38385476131666663)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(10, $i + 1).Value = $row[$i]
}
$row = @(-0.008780286409333371,0.74239403134999993,1.3116883506666666,0.93163655596166661,0.99582303529999994,1.0473839364000002,0.8639629678966666,0.74142034705000004,0.6768387450808333,0.74310114763333335,1.2205505145366664,1.3785349710500001,0.91431706577999994,0.52809710736666671)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(11, $i + 1).Value = $row[$i]
}
$row = @(0.23760426128833334,0.19504131647000003,0.27992919274183337,0.18352322803333335,0.28234097705,-0.10689199228333336,-0.16022735839833332,0.45700010125666662,0.95867359485333326,-0.2572886439,-0.67903045423333341,0.15027928414999997,0.033971047599999998,0.27846535920000004)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(12, $i + 1).Value = $row[$i]
}
$row = @(1.0133406782500001,1.1716971371166667,1.1802645971066665,0.76419523500999997,0.89041618510833331,0.63570949903333329,0.67641289383000003,0.73114867604333322,0.73554241765000006,0.83714145075000013,0.84667350447566669,0.99885304480999992,0.95713038591666677,1.3678465312100001)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(13, $i + 1).Value = $row[$i]
}
$row = @(0.029155728400000018,-0.33937920623333334,0.21371935165,0.5689332012566668,0.22426946611666668,0.086916582620000016,-0.12465734305000004,0.13000588694999998,0.47540883956333335,-0.11048061751666666,-0.46928690011999991,0.25380091815333333,0.27502195716666666,0.0068387690233333298)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(14, $i + 1).Value = $row[$i]
}
$row = @(-0.28976546560499999,0.090899667283333374,-0.029823568666666689,-0.24128871923833339,0.16270305365000004,0.29519310658333336,0.7099144587166667,0.52987478033333324,-0.14416405883166666,0.39379968725000003,0.055956515783333272,-0.13280236675000012,0.26363818210000001,0.16720295503333338)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(15, $i + 1).Value = $row[$i]
}
$row = @(-0.14353823563166668,0.078752278671666703,-0.14056636071833334,-0.29334329376999996,-0.10461769801666662,-0.11039479799333335,-0.021789377350000028,-0.2311939119333333,-0.13754617425666665,0.0016051055666667091,0.026470237073333402,-0.59587255072666667,-0.27081921410666671,0.09049800395999999)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(16, $i + 1).Value = $row[$i]
}
$row = @(-0.26428252388333334,-0.32754027049666667,0.0099633204999999836,0.16669843709999999,0.15666989588333335,0.0065458240633333581,-0.28722640417666667,-0.5967308204166667,-0.3852605231366667,-0.85414740671666667,0.044911147573333365,-0.57618666003666674,-0.36933440941666662,-1.0721190213333334)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(17, $i + 1).Value = $row[$i]
}
$row = @(-0.088234786499999995,-0.21576441113333339,0.27704523859999997,0.21672758729999997,0.65492175697666666,0.10536115491333331,-0.31220379706666668,-0.32990892126000004,-0.14528807345666661,-0.2604948765766667,0.28393472004999998,-0.8141420855,-1.3404253371500001,-0.50616982481666672)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(18, $i + 1).Value = $row[$i]
}
$row = @(-0.11004771014000002,0.25802007916999997,0.084481388566666676,0.5913919241333333,0.10677577316333331,0.035143676766666687,-0.49200280937999996,-0.57495112673333337,0.37292041346666671,-1.0118914308333333,0.11656647340000001,-0.36976234951999998,0.93813614749666663,0.32293528766666668)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(19, $i + 1).Value = $row[$i]
}
$row = @(0.64604896271333334,0.37667879403333332,0.48875589225999999,0.043550502656666663,0.22340044590000002,0.38103119526333334,0.23953632953333334,-0.16390435176666662,0.37251095180000005,0.49988569395999999,0.28509423863666661,-0.25848657701666666,-0.18794973721666666,0.1445339541166667)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(20, $i + 1).Value = $row[$i]
}
$row = @(0.3222654177,0.048137901943333339,0.012200184686666643,0.38833046339333332,0.40613343508333333,0.38494016820000004,0.12059365480000003,0.2885819411433333,0.23231048809999999,0.93301154236666661,-0.0027429540333332447,-0.034515976120000014,-0.1468798457833333,0.12039870609999995)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(21, $i + 1).Value = $row[$i]
}
$row = @(-0.29056115942999999,-0.14286762083333332,0.038161567833333299,-0.28471795532666666,-0.46860664937333335,0.48569643431999998,0.56957966127366666,0.20203681714333332,0.10425142066666668,-0.71095812829666671,0.18892250415333334,-0.024601536333333333,0.11423210706666666,0.17548517249999998)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(22, $i + 1).Value = $row[$i]
}
$row = @(-0.37722851204999996,-1.054850586975,-0.47104815304999992,0.167037502945,0.10044010610999998,0.08684765709999992,0.44412407394999998,0.15436374764999997,-0.34735497095000001,-0.52173046679999979,-0.6017670529500001,-0.21720577125000001,0.55088863370000007,-0.59876922725000015)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(23, $i + 1).Value = $row[$i]
}
$row = @(0.071012296318333307,1.1171305833333334,1.2171792086499997,0.10293346568333334,0.35420033723333333,0.078202584250000012,0.058993754671666676,0.059705275173333322,0.58967174394999999,-0.40164761666666671,1.3547201436,0.4053123532166667,0.68644499741666665,0.7089177011333333)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(24, $i + 1).Value = $row[$i]
}
$row = @(0.10484316466666668,0.42420217983166664,0.62754424063000003,-0.40002227009166663,-0.54338901281666674,0.26980713584000005,0.33678950336666669,0.12332444751333338,-0.61881769894000005,0.16087607899333334,0.37308547261666658,0.23232726104200002,0.36274531348666672,-0.22496011209000003)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(25, $i + 1).Value = $row[$i]
}
$row = @(0.21857637891999998,-0.014252640519999998,0.15925035007499999,-0.57583307854999999,-0.53816728975000006,-0.2946449458,0.20432715093000001,0.35117492919999993,-0.23915359410000003,-0.13649009412999999,0.72739652911999997,0.75536674713999996,0.75517603300000002,0.43202635257599997)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(26, $i + 1).Value = $row[$i]
}
$row = @(-0.0090032762666666377,0.14477088027999999,0.0049571659683333158,-0.95575277110000001,-0.8757968661333333,-0.51380207443300008,-0.24666518500333334,-0.2711417404,-0.54714860363333329,0.4929259644916667,-0.84795449898666675,-0.51317119542833334,-0.092873170833333379,0.10282881825000001)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(27, $i + 1).Value = $row[$i]
}
$row = @(-0.16419693229833332,-0.14535240291833332,0.42048883753166666,0.46209290828333338,0.2610584630833333,-0.13345141705333333,-0.34227676643333327,-0.1195930077983333,-0.13650355778333334,-0.39216482439666661,0.075505691383333351,-0.31352124298219997,-0.18004938609999999,-0.65304382198333333)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(28, $i + 1).Value = $row[$i]
}
$row = @(0.40575048191666668,0.52861184906666658,0.22897109874999999,-0.26672506517333339,0.084619309133333265,0.3225371527333335,0.47615029321666663,0.65929490320999995,0.5251336253833333,1.0500000755499999,0.095617957266666664,0.24633849840333333,-0.98245840251666672,0.35693244628333343)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(29, $i + 1).Value = $row[$i]
}
$row = @(-0.19282722453333331,-0.010840832816666646,0.30994059416666664,0.12761237948333334,-0.1342137442766666,0.17478039148333332,0.076572034090000005,-0.034368030215000005,0.047176334083333354,-0.21340750930000002,-0.45169441438333341,-0.50708977350999995,0.054993370684999987,0.28344589966666667)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(30, $i + 1).Value = $row[$i]
}
$row = @(-0.53084684260000004,-0.08111045440833331,-0.012817388216666649,-0.34240139503333333,-1.0984878487833334,-0.81726936004999995,-0.32354331586666663,-0.63182363183666668,-0.67982109755333331,-0.79024459480000009,-0.60944920487999998,-0.52049932698333334,-0.27796738529999998,-0.27601842881666672)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(31, $i + 1).Value = $row[$i]
}
$row = @(0.43516055118999997,0.60486066879,0.13971840147666667,0.17075054015333335,0.45389919980333332,0.011662660950000008,0.097388882403333363,0.24860731887333337,0.51803874802999994,0.021382271456666671,-0.26444643716666666,0.37297881628333335,-0.16439958338333338,0.4046217375233333)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(32, $i + 1).Value = $row[$i]
}
$row = @(0.5227404881833333,0.76928472041666662,0.67858297254999989,0.20565971838333336,0.55896325645333322,0.32620372698333333,-0.15897013125166665,-0.59841555933666668,-0.16859921024999999,-1.4495397248000002,0.80525976830000001,0.88534372598333344,0.96416393864333338,0.33995330447333327)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(33, $i + 1).Value = $row[$i]
}
$row = @(0.48069648450000002,-0.67936644646666666,0.16356151522000001,0.85133621943333337,0.55168261527666673,0.32272039474999997,0.7228148386333334,0.46695444906666667,0.88563983406666669,1.0024038527999999,0.34231524932333335,-0.10494660499666665,0.80232802263333336,0.007469059499999986)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(34, $i + 1).Value = $row[$i]
}
$row = @(0.38976165965000004,0.86641512923333319,0.88408703932333332,0.19997032281666668,0.44602612487333337,0.9047739664929999,1.0570409516133332,0.44821938664333333,0.52492414330000003,0.73198449674999999,0.872774991035,0.47008789784499999,0.4093129386183334,0.66443219755000005)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(35, $i + 1).Value = $row[$i]
}
$row = @(-0.59191742235000011,-0.42563440268333336,-0.38669340585,0.043613281129999987,0.23612410184933333,-0.26576967823333336,-0.53255249351666678,-0.14112532831666669,0.27099792215000001,0.128090652423,-0.86243856381666661,-0.73290970948166678,-0.80998644841666667,-0.9987095444716666)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(36, $i + 1).Value = $row[$i]
}
$row = @(-0.079296892250000028,-0.26554660297166666,0.050218083289999986,0.48702313459166668,-0.27440186274000011,0.51703115266666666,0.52838426366500002,0.24180795766666668,0.57581442824666673,1.2240870836833335,0.026188398369999998,0.7559955192333333,0.29372968554000006,0.95347662811666667)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(37, $i + 1).Value = $row[$i]
}
$row = @(-0.30822447019999999,-0.36457650326666663,-0.20133795677166669,0.25754077167166667,0.42439166029999997,-0.10352517669016668,-0.10308745100333327,0.092796559916666632,-0.26948277283396666,0.57043169846666664,0.11187272138833335,-0.15370623980000003,-0.29276354888999995,-1.01261724536)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(38, $i + 1).Value = $row[$i]
}
$row = @(-0.015952305190000016,0.094068659549999989,-0.083657220416666678,0.77000307288333325,0.75897578340000005,0.27884429630000007,-0.052927063216666626,0.4322770917,0.22665833173166666,-0.52389736335000003,-0.31043187090000002,0.31903231523666664,1.3113785312666666,0.016093494183333334)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(39, $i + 1).Value = $row[$i]
}
$row = @(0.36264404996666666,0.66912919101666668,0.63446635352333347,-0.64798807383333323,-0.50418314438333334,0.034367229983333358,0.35491979816666674,0.23503923550333333,-0.095598167104999993,-0.19290804825000007,0.99262397550000014,0.65631808275666659,0.25619537553333338,0.31100248281666676)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(40, $i + 1).Value = $row[$i]
}
$row = @(0.68179065840666675,0.90013687649666663,0.99629175808333326,0.11951657348333333,-0.21315500751666666,0.39620110896666672,0.93065795703166665,0.61998084655666663,0.46674710580000001,1.0760125878266666,0.99040308406666666,0.80600059545000002,1.0614589461999999,1.2041940128366666)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(41, $i + 1).Value = $row[$i]
}
$row = @(0.21467334996666665,0.22188899451166663,-0.18298943671500004,-0.0058921143500000328,-0.39274086753333332,-0.06772546829999998,-0.11940955946666666,0.23519633063000001,0.15447596626666668,0.80037809185000008,-0.060285529216666689,0.7485319109666666,-0.16452812328333333,0.38171453296666669)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(42, $i + 1).Value = $row[$i]
}
$row = @(-0.20539374558333334,-0.56134917258566663,-0.15234341770000001,0.18055077061666669,-0.12422792028333335,0.43199212126666664,0.28222207263333327,0.2594102242752,0.11651138188333343,1.3117781174666665,-0.49987177821666662,-0.41120418219666666,-0.82183717733999995,-0.21255795461666671)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(43, $i + 1).Value = $row[$i]
}
$row = @(-0.15394203779233334,0.245653184425,0.41404490253333337,-0.15111926351999999,-0.11432460240449999,0.46271834861666672,0.8432669500833333,0.37598362513333333,0.13797137045999999,0.83987219717666672,0.31909237840399995,0.24012425828566666,0.88048015582333328,-0.35408785891666672)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(44, $i + 1).Value = $row[$i]
}
$row = @(0.18401484926666667,0.71069887840000001,1.0865463314733332,0.51840498217666653,0.11435656117500009,0.31817600443999999,0.77842012738666677,-0.028052262959333343,-0.18158266191,-0.04973520569333334,0.92753450578999996,0.74105106708333335,0.64549219276666658,0.74516944586666667)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(45, $i + 1).Value = $row[$i]
}
$row = @(-0.26689443811499997,-0.54488971445000001,0.12086491939999994,-0.22940440835000009,-0.17579611114999999,-0.66041109345000004,-0.29548238000000004,-0.36232935654999998,-0.28072121915000003,-0.27172762644999993,-0.085885630749999969,0.11025550239999998,0.2682410348,0.38784159304999993)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(46, $i + 1).Value = $row[$i]
}
$row = @(-0.14381807821000001,-0.18104740728333327,-0.74068978794999996,-0.21944982681266667,-0.073709379796666641,-0.60236857593333326,-0.75543666661666664,-0.16596575768666669,-0.26959552923333341,-0.71806960630000005,-0.68031695271666681,-0.61568571700000008,0.32735701903333336,-0.058076726896999942)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(47, $i + 1).Value = $row[$i]
}
$row = @(0.312435472775,0.54161556929999999,0.48000072199999999,0.050000955547499976,0.43351274662000006,0.28724576245,0.69635890899999997,0.64256019208999993,0.06466816265,1.5581649125500001,0.82639124725000002,-0.093605462149999985,0.17822328260000006,1.1797739003000001)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(48, $i + 1).Value = $row[$i]
}
$row = @(0.74734466240499997,1.2341783366499999,1.0675437705999999,-0.028411920224999943,-0.25606632201499996,-0.056371029010000018,-0.15997358736049999,-0.13518493195,0.34696181728499997,0.39087414624999994,1.3647107410000001,0.78117012369999994,0.91885792170000002,0.99794760690000006)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(49, $i + 1).Value = $row[$i]
}
$row = @(0.1547678319583333,0.15408608190966669,-0.14162731094333336,0.36181857151666669,-0.12576519458333332,0.16136178383333333,0.49447943797666666,0.36916307965,0.34368686159166661,0.33780602268333332,-0.37864099104333337,-0.13191306091499999,-0.20642586136666666,-0.13872520213666667)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(50, $i + 1).Value = $row[$i]
}
$row = @(-0.40399291024999995,-1.1309873159700001,-1.4385633276500003,-1.2130543523000001,-1.8274642631,-0.76997308745000004,-0.20187996105000006,-0.19399287195000003,-0.99824027179999997,0.91696570410000011,-1.4894017379500002,-1.2962486646,-1.4480606497999999,-1.8229347162500003)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(51, $i + 1).Value = $row[$i]
}
$row = @(-0.3070672176,-0.16743448005500006,0.28639530935000002,0.11390298285000006,0.086239355045000021,0.17319648134999993,-0.17734284075000001,0.17968162034999999,0.21792126275500004,0.083064706719999992,0.14325645308000001,-0.32314184099999999,-0.75208651944999994,-0.041033311605000011)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(52, $i + 1).Value = $row[$i]
}
$row = @(-0.9890844929,-0.64983048141599997,-0.35770830791500002,0.24515135911500002,0.56471008695000002,0.44038670785,0.11992243570050004,0.12742385479499999,-0.13914139639999998,0.53612783740000003,1.1256037351500001,0.26828021859000001,-0.034885794349999943,-0.73703056238999998)
for ($i = 0; $i -lt $row.Length; $i++) {
  $ws.Cells.Item(53, $i + 1).Value = $row[$i]
}

# --- Column width adjustments ---
$ws.Columns(2).ColumnWidth = 14.666666666666666
$ws.Columns(3).ColumnWidth = 13.833333333333334
$ws.Columns(5).ColumnWidth = 14.833333333333334
$ws.Columns(7).ColumnWidth = 13.666666666666666
$ws.Columns(9).ColumnWidth = 12.833333333333334
$ws.Columns(10).ColumnWidth = 13.833333333333334
$ws.Columns(11).ColumnWidth = 14.666666666666666
$ws.Columns(14).ColumnWidth = 13.833333333333334
